$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Phút hành chính" column (column K) entirely, shifting
# everything to its right one column to the left.
$ws.Range("K1").EntireColumn.Delete()

# Update the selected view to roughly match the post-edit state.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("M8").Select()
